$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'basketball leg sleeve youth boys'
$ws.Range('A2').Value = 'youth basketball leg sleeves boys'
$ws.Range('A3').Value = 'basketball pants with knee pads'
$ws.Range('A4').Value = 'mens compression pants'
$ws.Range('A5').Value = 'compression knee pad'
$ws.Range('A6').Value = 'basketball leggings'
$ws.Range('A7').Value = 'knee compression pads'
$ws.Range('A8').Value = 'compression pants men'
$ws.Range('A9').Value = 'knee pads hex'
$ws.Range('A10').Value = 'basketball knee pads youth boys'
$ws.Range('A11').Value = 'mens capri pants'
$ws.Range('A12').Value = 'volleyball knee'
$ws.Range('A13').Value = 'workout hand pads'
$ws.Range('A14').Value = 'black athletic leggings'
$ws.Range('A15').Value = 'man pads'
$ws.Range('A16').Value = 'compression tights'
$ws.Range('A17').Value = 'tactical pants with knee pads'
$ws.Range('A18').Value = 'mens pants'
$ws.Range('A19').Value = 'xxl knee pads'
$ws.Range('A20').Value = 'knee compression for men'
$ws.Range('A21').Value = 'athletic leggings capri'
$ws.Range('A22').Value = 'sliding pads for softball'
$ws.Range('A23').Value = 'knee pads adult'
$ws.Range('A24').Value = 'softball knee pads'
$ws.Range('A25').Value = 'basketball knee pads adult'
$ws.Range('A26').Value = 'youth basketball gear'
$ws.Range('A27').Value = 'boys leggings'
$ws.Range('A28').Value = 'compression knee pads'
$ws.Range('A29').Value = 'black basketball leggings for men'
$ws.Range('A30').Value = 'compression pants with knee pads boys'
$ws.Range('A31').Value = 'hex knee pads basketball youth'
$ws.Range('A32').Value = 'anti strip clothing for men'
$ws.Range('A33').Value = 'hex protective knee pads'
$ws.Range('A34').Value = 'men leggings'
$ws.Range('A35').Value = 'weightlifting pads'
$ws.Range('A36').Value = 'compression tights boys'
$ws.Range('A37').Value = 'compression tights for boys'
$ws.Range('A38').Value = 'bjj tights for men'
$ws.Range('A39').Value = 'mens football pants'
$ws.Range('A40').Value = 'advanced squat pad'
$ws.Range('A41').Value = 'knee pads football'
$ws.Range('A42').Value = 'knee protector pain'
$ws.Range('A43').Value = 'knee support leggings'
$ws.Range('A44').Value = 'compression tights for youth'
$ws.Range('A45').Value = 'mens tall pants'
$ws.Range('A46').Value = 'men tights'
$ws.Range('A47').Value = 'wrestling gear'
$ws.Range('A48').Value = 'little boys athletic leggings'
$ws.Range('A49').Value = 'workout leggings for men pack'
$ws.Range('A50').Value = 'baseball gear for boys'
$ws.Range('A51').Value = 'black baseball pants youth xl'
$ws.Range('A52').Value = 'wrestling knee pad'
$ws.Range('A53').Value = 'volleyball knee pads xxl mens'
$ws.Range('A54').Value = 'knee running'
$ws.Range('A55').Value = 'knee pads for gym'
$ws.Range('A56').Value = 'knee for running'
$ws.Range('A57').Value = 'athletic compression leggings'
$ws.Range('A58').Value = 'weightlifting pad'
$ws.Range('A59').Value = 'foam knee pad'
$ws.Range('A60').Value = 'compression tight'
$ws.Range('A61').Value = 'knee support volleyball'
$ws.Range('A62').Value = 'basketball equipment'
$ws.Range('A63').Value = 'fitness squat pad'
$ws.Range('A64').Value = 'knee swelling'
$ws.Range('A65').Value = 'softball pants mens black'
$ws.Range('A66').Value = 'sport pants for men'
$ws.Range('A67').Value = 'basketball equipment pads'
$ws.Range('A68').Value = 'knee length leggings'
$ws.Range('A69').Value = 'volleyball knee pads girls youth'
$ws.Range('A70').Value = 'reduce swelling after surgery'
$ws.Range('A71').Value = 'soccer protection'
$ws.Range('A72').Value = 'knee support for basketball men'
$ws.Range('A73').Value = 'leg compression for men'
$ws.Range('A74').Value = 'knee support for basketball'
$ws.Range('A75').Value = 'soccer clothes'
$ws.Range('A76').Value = 'lacrosse equipment'
$ws.Range('A77').Value = 'capri athletic pants'
$ws.Range('A78').Value = 'basketball knee pads youth pair'
$ws.Range('A79').Value = 'mens leggings compression winter'
$ws.Range('A80').Value = 'ski compression pants men'
$ws.Range('A81').Value = 'leggings pants'
$ws.Range('A82').Value = 'youth volleyball knee pads'
$ws.Range('A83').Value = 'nike leggings basketball men'
$ws.Range('A84').Value = 'underarmour tights'
$ws.Range('A85').Value = 'yourh basketball knee pads'
$ws.Range('A86').Value = 'compression leggings knee pads'
$ws.Range('A87').Value = 'basketball tights with knee pads'
$ws.Range('A88').Value = 'leggings for men with knee pads'
$ws.Range('A89').Value = 'basketball knee pads men'
$ws.Range('A90').Value = 'mens compression capri pants'
$ws.Range('A91').Value = 'lacrosse leggings'
$ws.Range('A92').Value = 'mens compression leggings capri'
$ws.Range('A93').Value = 'basketball hip pads'
$ws.Range('A94').Value = 'male compression tights'
$ws.Range('A95').Value = 'compression tights for girls'
$ws.Range('A96').Value = 'basketball knee pads'
$ws.Range('A97').Value = 'mens compression pants capri'
$ws.Range('A98').Value = 'compression mens pants'
$ws.Range('A99').Value = 'basketball leggings for men'
$ws.Range('A100').Value = 'youth capri compression tights'
